# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same 7-column fund-holding layout) to
#    create a new sheet positioned right before "总计", rename it "2022-Q1",
#    and overwrite its data with the new quarter's fund holdings.
# 2. Update the "总计" (summary) sheet: insert a new top data row for
#    "2022-Q1" and push the existing 2021-Q4 / 2021-Q2 rows down.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet  = $wb.Worksheets.Item("总计")

# --- 1. Create "2022-Q1" sheet (copy of "2021-Q4", placed right before 总计) ---
$sourceSheet.Copy($totalSheet, $null)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Force the data columns (B:G) to stay text so numeric-looking values like
# "009225" / "1.84" are not auto-coerced to numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "009225"
$newSheet.Range("C2").Value = "天弘中证中美互联网指数（QDII）A"
$newSheet.Range("D2").Value = "1.84"
$newSheet.Range("E2").Value = "94.90"
$newSheet.Range("F2").Value = "6.30"
$newSheet.Range("G2").Value = "0.1159"
$newSheet.Range("H2").Value = 6

$newSheet.Range("B3").Value = "009226"
$newSheet.Range("C3").Value = "天弘中证中美互联网指数（QDII）C"
$newSheet.Range("D3").Value = "0.59"
$newSheet.Range("E3").Value = "94.90"
$newSheet.Range("F3").Value = "6.30"
$newSheet.Range("G3").Value = "0.0372"
$newSheet.Range("H3").Value = 6

# --- 2. Update "总计" sheet with the new 2022-Q1 row on top ---
# Re-fetch by name: inserting/copying a sheet shifts tab positions, and the
# earlier $totalSheet handle can end up pointing at the wrong tab.
$totalSheet = $wb.Worksheets.Item("总计")

# Carry the row-1-style (centered/bordered) formatting down onto the new A4
# cell before writing into it (previously-unused cell).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.12

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.05

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.15

# Restore the originally-active tab: copying a sheet makes the new copy the
# active one, but "2021-Q2" was selected before this edit.
$wb.Worksheets.Item("2021-Q2").Activate()
